$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "Price" column (D) values: protect as Text so Excel does not
# --- reinterpret numeric-looking strings (e.g. "186.10" -> 186.1, "0.0000251" -> 2.51E-05)
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '65.543.12'
$ws.Range("D2").NumberFormat = "General"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.487.15'
$ws.Range("D3").NumberFormat = "General"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").NumberFormat = "General"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '565.34'
$ws.Range("D5").NumberFormat = "General"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '186.10'
$ws.Range("D6").NumberFormat = "General"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.485.18'
$ws.Range("D7").NumberFormat = "General"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.600'
$ws.Range("D8").NumberFormat = "General"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.650'
$ws.Range("D10").NumberFormat = "General"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.141'
$ws.Range("D11").NumberFormat = "General"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '51.76'
$ws.Range("D12").NumberFormat = "General"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000251'
$ws.Range("D13").NumberFormat = "General"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '9.52'
$ws.Range("D14").NumberFormat = "General"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.043.25'
$ws.Range("D15").NumberFormat = "General"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.489.81'
$ws.Range("D17").NumberFormat = "General"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '17.88'
$ws.Range("D18").NumberFormat = "General"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '65.293.76'
$ws.Range("D19").NumberFormat = "General"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.81'
$ws.Range("D20").NumberFormat = "General"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '382.68'
$ws.Range("D22").NumberFormat = "General"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.17'
$ws.Range("D23").NumberFormat = "General"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '83.43'
$ws.Range("D24").NumberFormat = "General"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '10.75'
$ws.Range("D25").NumberFormat = "General"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.82'
$ws.Range("D26").NumberFormat = "General"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '12.00'
$ws.Range("D28").NumberFormat = "General"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '3.41'
$ws.Range("D29").NumberFormat = "General"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.62'
$ws.Range("D30").NumberFormat = "General"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '30.25'
$ws.Range("D31").NumberFormat = "General"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '11.86'
$ws.Range("D33").NumberFormat = "General"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '606.91'
$ws.Range("D34").NumberFormat = "General"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '62.46'
$ws.Range("D35").NumberFormat = "General"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '40.45'
$ws.Range("D37").NumberFormat = "General"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.388'
$ws.Range("D39").NumberFormat = "General"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0₃0730'
$ws.Range("D40").NumberFormat = "General"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.934.91'
$ws.Range("D43").NumberFormat = "General"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.73'
$ws.Range("D44").NumberFormat = "General"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.41'
$ws.Range("D45").NumberFormat = "General"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0394'
$ws.Range("D46").NumberFormat = "General"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.03'
$ws.Range("D47").NumberFormat = "General"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.127'
$ws.Range("D48").NumberFormat = "General"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.21'
$ws.Range("D50").NumberFormat = "General"

# --- Update remaining cells (Coin name, Link, Volume(1h) percentage text)
$ws.Range("E2").Value = '  -5.27%  '
$ws.Range("E3").Value = '  -6.41%  '
$ws.Range("E4").Value = '  +0.28%  '
$ws.Range("E5").Value = '  -7.69%  '
$ws.Range("E6").Value = '  -1.89%  '
$ws.Range("E7").Value = '  -6.24%  '
$ws.Range("E8").Value = '  -5.83%  '
$ws.Range("E9").Value = '  -0.11%  '
$ws.Range("E10").Value = '  -9.99%  '
$ws.Range("E11").Value = '  -12.20%  '
$ws.Range("E12").Value = '  -12.18%  '
$ws.Range("E13").Value = '  -13.72%  '
$ws.Range("E14").Value = '  -10.94%  '
$ws.Range("E15").Value = '  -6.21%  '
$ws.Range("E16").Value = '  -1.39%  '
$ws.Range("E17").Value = '  -6.12%  '
$ws.Range("E18").Value = '  -7.63%  '
$ws.Range("E19").Value = '  -5.28%  '
$ws.Range("E20").Value = '  -8.86%  '
$ws.Range("E21").Value = '  -9.74%  '
$ws.Range("E22").Value = '  -7.24%  '
$ws.Range("E23").Value = '  -9.41%  '
$ws.Range("E24").Value = '  -6.69%  '
$ws.Range("E25").Value = '  -1.46%  '
$ws.Range("E26").Value = '  -7.77%  '
$ws.Range("E27").Value = '  -0.49%  '
$ws.Range("E28").Value = '  -6.73%  '
$ws.Range("E29").Value = '  -10.49%  '
$ws.Range("E30").Value = '  -10.98%  '
$ws.Range("E31").Value = '  -8.82%  '
$ws.Range("E32").Value = '  -8.20%  '
$ws.Range("B33").Value = 'Cosmos'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("E33").Value = '  -7.14%  '
$ws.Range("B34").Value = 'Bittensor'
$ws.Range("C34").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("E34").Value = '  -5.02%  '
$ws.Range("E35").Value = '  -4.86%  '
$ws.Range("E36").Value = '  -10.31%  '
$ws.Range("E37").Value = '  -12.71%  '
$ws.Range("E38").Value = '  +0.07%  '
$ws.Range("E39").Value = '  -6.12%  '
$ws.Range("E40").Value = '  -12.29%  '
$ws.Range("E41").Value = '  +0.07%  '
$ws.Range("E42").Value = '  -8.74%  '
$ws.Range("E43").Value = '  +2.31%  '
$ws.Range("E44").Value = '  -10.33%  '
$ws.Range("E45").Value = '  -8.37%  '
$ws.Range("E46").Value = '  -11.46%  '
$ws.Range("E47").Value = '  -1.62%  '
$ws.Range("E48").Value = '  -8.95%  '
$ws.Range("E49").Value = '  -3.65%  '
$ws.Range("E50").Value = '  -10.20%  '
$ws.Range("E51").Value = '  -10.23%  '
